# Apply the "update hotel reviews data" edit:
#  - hotel_info: fill English_Reviews_num / Local_Rank / Total_Reviews_num on row 2
#  - review_info: append the 6 newly scraped reviews (rows 2-7)
$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force text storage so numeric-/date-looking strings ("11", "08/08/2018",
    # "June 2018", ...) are written as shared strings, not coerced into
    # numbers/dates by Excels usual input parsing.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- hotel_info sheet: fill in review-count / rank fields on row 2 ---
$ws1 = $wb.Worksheets.Item("hotel_info")
Set-TextValue $ws1.Range("G2") '11'
Set-TextValue $ws1.Range("H2") '7'
Set-TextValue $ws1.Range("I2") '11'

# --- review_info sheet: add the 6 scraped review rows ---
$ws2 = $wb.Worksheets.Item("review_info")

# Row 2
$ws2.Range("A2").Value = 66432
$ws2.Range("D2").Value = 1
Set-TextValue $ws2.Range("E2") '08/08/2018'
Set-TextValue $ws2.Range("F2") 'https://www.tripadvisor.com/ShowUserReviews-g56723-d12553156-r588580962-Courtyard_Houston_Sugar_Land_Lake_Pointe-Sugar_Land_Texas.html'
Set-TextValue $ws2.Range("G2") '56723'
Set-TextValue $ws2.Range("H2") '12553156'
Set-TextValue $ws2.Range("I2") '588580962'
Set-TextValue $ws2.Range("J2") '06/18/2018'
Set-TextValue $ws2.Range("K2") 'Great, Nearly New Hotel'
Set-TextValue $ws2.Range("L2") 'We were attending a concert at Smart Financial Center and didn''t want to drive the 70 miles home late at night, so stayed here.  This Courtyard by Marriott is 10 minutes from this venue.  Room was exquisitely clean, and everything worked well.  All the staff, front desk Bistro Dining and Housekeeping were helpful, efficient and friendly.  They genuinely appreciated our business.  Exercise facility on site, and location is near a tree-lined walking trail, lake.  Great to have the Bistro on site, if you don''t want to drive to another restaurant.  We had dinner, which was excellent.  Many other restaurants over all price ranges are nearby.'
$ws2.Range("M2").Value = 5
Set-TextValue $ws2.Range("N2") 'June 2018'
Set-TextValue $ws2.Range("O2") ' traveled as a couple'
$ws2.Range("R2").Value = 5
$ws2.Range("S2").Value = 5
$ws2.Range("U2").Value = 5
$ws2.Range("V2").Value = 0
Set-TextValue $ws2.Range("Y2") 'We were attending a concert at Smart Financial Center and didn''t want to drive the 70 miles home late at night, so stayed here.  This Courtyard by Marriott is 10 minutes from this venue.  Room was exquisitely clean, and everything worked well.  All the staff, front desk Bistro Dining and Housekeeping were helpful, efficient and friendly.  They genuinely appreciated our business.  Exercise facility on site, and location is near a tree-lined walking trail, lake.  Great to have the Bistro on site, if you don''t want to drive to another restaurant.  We had dinner, which was excellent.  Many other restaurants over all price ranges are nearby.'

# Row 3
$ws2.Range("A3").Value = 66432
$ws2.Range("D3").Value = 2
Set-TextValue $ws2.Range("E3") '08/08/2018'
Set-TextValue $ws2.Range("F3") 'https://www.tripadvisor.com/ShowUserReviews-g56723-d12553156-r582935819-Courtyard_Houston_Sugar_Land_Lake_Pointe-Sugar_Land_Texas.html'
Set-TextValue $ws2.Range("G3") '56723'
Set-TextValue $ws2.Range("H3") '12553156'
Set-TextValue $ws2.Range("I3") '582935819'
Set-TextValue $ws2.Range("J3") '05/26/2018'
Set-TextValue $ws2.Range("K3") 'Like to run along the lake, open the door'
Set-TextValue $ws2.Range("L3") 'Sittting on a beautiful lake with shaded walking/jogging trails and 5 minutes from shopping, movies, restaurants would be enough but being able to sit out and eat on the Lake day or night is a bonus. I did my paperwork outdoors. Sitting and eating areas are shaded. The outdoor pool is on the lake also. In the evenings there is a fire pit. Wish they sold s’more supplies. Be sure to package and take your own. The double shower with the non slip shower was beautiful and of course the BED the awesome Marriott bed. Food offerings were good. The oatmeal bowl and the Brioche breakfast sandwich were my choices.MoreShow less'
$ws2.Range("M3").Value = 5
Set-TextValue $ws2.Range("N3") 'May 2018'
Set-TextValue $ws2.Range("O3") ' traveled on business'
$ws2.Range("V3").Value = 0
Set-TextValue $ws2.Range("Y3") 'Sittting on a beautiful lake with shaded walking/jogging trails and 5 minutes from shopping, movies, restaurants would be enough but being able to sit out and eat on the Lake day or night is a bonus. I did my paperwork outdoors. Sitting and eating areas are shaded. The outdoor pool is on the lake also. In the evenings there is a fire pit. Wish they sold s’more supplies. Be sure to package and take your own. The double shower with the non slip shower was beautiful and of course the BED the awesome Marriott bed. Food offerings were good. The oatmeal bowl and the Brioche breakfast sandwich were my choices.More'

# Row 4
$ws2.Range("A4").Value = 66432
$ws2.Range("D4").Value = 3
Set-TextValue $ws2.Range("E4") '08/08/2018'
Set-TextValue $ws2.Range("F4") 'https://www.tripadvisor.com/ShowUserReviews-g56723-d12553156-r573527011-Courtyard_Houston_Sugar_Land_Lake_Pointe-Sugar_Land_Texas.html'
Set-TextValue $ws2.Range("G4") '56723'
Set-TextValue $ws2.Range("H4") '12553156'
Set-TextValue $ws2.Range("I4") '573527011'
Set-TextValue $ws2.Range("J4") '04/15/2018'
Set-TextValue $ws2.Range("K4") 'Nice, new, easy to get to, comfy beds'
Set-TextValue $ws2.Range("L4") 'Great price for a clean, comfy room in a great location. Room was very nice. The bed was very comfortable. The king room had a love seat and ottoman that could have slept another small adult or child. Bathroom had a large walk-in shower and plenty of counter space for us women. There is a Starbucks/bar/snack area downstairs (nothing free) wine and beer was quite expensive. Would definitely stay again.'
$ws2.Range("M4").Value = 4
Set-TextValue $ws2.Range("N4") 'April 2018'
Set-TextValue $ws2.Range("O4") ' traveled as a couple'
$ws2.Range("V4").Value = 0
Set-TextValue $ws2.Range("Y4") 'Great price for a clean, comfy room in a great location. Room was very nice. The bed was very comfortable. The king room had a love seat and ottoman that could have slept another small adult or child. Bathroom had a large walk-in shower and plenty of counter space for us women. There is a Starbucks/bar/snack area downstairs (nothing free) wine and beer was quite expensive. Would definitely stay again.'

# Row 5
$ws2.Range("A5").Value = 66432
$ws2.Range("D5").Value = 4
Set-TextValue $ws2.Range("E5") '08/08/2018'
Set-TextValue $ws2.Range("F5") 'https://www.tripadvisor.com/ShowUserReviews-g56723-d12553156-r536658353-Courtyard_Houston_Sugar_Land_Lake_Pointe-Sugar_Land_Texas.html'
Set-TextValue $ws2.Range("G5") '56723'
Set-TextValue $ws2.Range("H5") '12553156'
Set-TextValue $ws2.Range("I5") '536658353'
Set-TextValue $ws2.Range("J5") '10/28/2017'
Set-TextValue $ws2.Range("K5") 'Mixed review'
Set-TextValue $ws2.Range("L5") 'The staff was friendly and our room was clean. The vent system forces you to hear other people’s conversations and the plumbing alerted you to when your neighbors flushed the toilet. The bed was very uncomfortable. '
$ws2.Range("M5").Value = 4
Set-TextValue $ws2.Range("N5") 'October 2017'
Set-TextValue $ws2.Range("O5") ' traveled as a couple'
$ws2.Range("V5").Value = 0
Set-TextValue $ws2.Range("Y5") 'The staff was friendly and our room was clean. The vent system forces you to hear other people’s conversations and the plumbing alerted you to when your neighbors flushed the toilet. The bed was very uncomfortable. '

# Row 6
$ws2.Range("A6").Value = 66432
$ws2.Range("D6").Value = 5
Set-TextValue $ws2.Range("E6") '08/08/2018'
Set-TextValue $ws2.Range("F6") 'https://www.tripadvisor.com/ShowUserReviews-g56723-d12553156-r531678791-Courtyard_Houston_Sugar_Land_Lake_Pointe-Sugar_Land_Texas.html'
Set-TextValue $ws2.Range("G6") '56723'
Set-TextValue $ws2.Range("H6") '12553156'
Set-TextValue $ws2.Range("I6") '531678791'
Set-TextValue $ws2.Range("J6") '10/10/2017'
Set-TextValue $ws2.Range("K6") 'Brand New!'
Set-TextValue $ws2.Range("L6") 'This hotel was not yet open 1 week when I stayed here the first week in October. I recommend this over the Marriott across the highway in citywalk because parking is easier and you have a lovely area outside to sit fireside in front of the lake. The rooms are large and the staff is friendly. It also has a kitchen for all three meals and they serve Starbucks.'
$ws2.Range("M6").Value = 4
Set-TextValue $ws2.Range("N6") 'October 2017'
Set-TextValue $ws2.Range("O6") ' traveled on business'
$ws2.Range("R6").Value = 5
$ws2.Range("S6").Value = 5
$ws2.Range("U6").Value = 5
$ws2.Range("V6").Value = 0
Set-TextValue $ws2.Range("Y6") 'This hotel was not yet open 1 week when I stayed here the first week in October. I recommend this over the Marriott across the highway in citywalk because parking is easier and you have a lovely area outside to sit fireside in front of the lake. The rooms are large and the staff is friendly. It also has a kitchen for all three meals and they serve Starbucks.'

# Row 7
$ws2.Range("A7").Value = 66432
$ws2.Range("D7").Value = 6
Set-TextValue $ws2.Range("E7") '08/08/2018'
Set-TextValue $ws2.Range("F7") 'https://www.tripadvisor.com/ShowUserReviews-g56723-d12553156-r528798547-Courtyard_Houston_Sugar_Land_Lake_Pointe-Sugar_Land_Texas.html'
Set-TextValue $ws2.Range("G7") '56723'
Set-TextValue $ws2.Range("H7") '12553156'
Set-TextValue $ws2.Range("I7") '528798547'
Set-TextValue $ws2.Range("J7") '10/01/2017'
Set-TextValue $ws2.Range("K7") 'Wonderful hotel in great location with great front desk service'
Set-TextValue $ws2.Range("L7") 'Teejay, lola and daja were wonderful at front desk answering a few questions about the amenities nearby. They showed me where the nice places to eat around the hotel  and had staff bring extra towels in a prompt fashion.fatimah was awesome as well'
$ws2.Range("M7").Value = 5
Set-TextValue $ws2.Range("N7") 'September 2017'
Set-TextValue $ws2.Range("O7") ' traveled solo'
$ws2.Range("Q7").Value = 5
$ws2.Range("S7").Value = 5
$ws2.Range("U7").Value = 5
$ws2.Range("V7").Value = 0
Set-TextValue $ws2.Range("Y7") 'Teejay, lola and daja were wonderful at front desk answering a few questions about the amenities nearby. They showed me where the nice places to eat around the hotel  and had staff bring extra towels in a prompt fashion.fatimah was awesome as well'

